$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = "2017-12-31 00:00:00"

$ws.Range("O2").Value = 64810282.72
$ws.Range("P2").Value = 828.8369748534
$ws.Range("Q2").Value = 538153026.38
$ws.Range("R2").Value = 6882.2586119555
$ws.Range("S2").Value = 56896663.18
$ws.Range("T2").Value = 727.6323479887
$ws.Range("U2").Value = -55438832.85
$ws.Range("V2").Value = -708.988644708

$ws.Range("W2").Value = 35479.41
$ws.Range("X2").Value = 0.4537342783

$ws.Range("Y2").Value = 55724312.26
$ws.Range("Z2").Value = 712.6395451614
$ws.Range("AA2").Value = -14378612.23
$ws.Range("AB2").Value = -183.8832506686
$ws.Range("AC2").Value = -7819424.65

$ws.Range("AD2").ClearContents()
